# Auto-generated Excel COM-interop script
# Updates FFXIV leve market-price data (columns H-N) across all 8 job sheets
# to match the scheduled-runner refresh described in the commit message.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 129.8
$ws.Range("I11").Value = 129.8
$ws.Range("K11").Value = 129.8
$ws.Range("M11").Value = 10.19999999999999
$ws.Range("H98").Value = 1480.3
$ws.Range("J98").Value = 1233
$ws.Range("L98").Value = 1233
$ws.Range("N98").Value = -4229
$ws.Range("H115").Value = 119.5
$ws.Range("I115").Value = 119.5
$ws.Range("K115").Value = 358.5
$ws.Range("M115").Value = 1208.5
$ws.Range("H116").Value = 4487.2
$ws.Range("I116").Value = 3339
$ws.Range("K116").Value = 3339
$ws.Range("M116").Value = 103
$ws.Range("H122").Value = 1480.3
$ws.Range("J122").Value = 1233
$ws.Range("L122").Value = 3699
$ws.Range("N122").Value = -8599
$ws.Range("H127").Value = 1738.4
$ws.Range("I127").Value = 1423
$ws.Range("K127").Value = 4269
$ws.Range("M127").Value = 691
$ws.Range("H129").Value = 2339.6365
$ws.Range("I129").Value = 534
$ws.Range("K129").Value = 1602
$ws.Range("M129").Value = 3398
$ws.Range("H137").Value = 1505.3334
$ws.Range("I137").Value = 1505.3334
$ws.Range("K137").Value = 4516.0002
$ws.Range("M137").Value = -1966.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 489.6
$ws.Range("I2").Value = 489.6
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 489.6
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -376.6
$ws.Range("N2").ClearContents()
$ws.Range("H5").Value = 310
$ws.Range("I5").Value = 310
$ws.Range("K5").Value = 310
$ws.Range("M5").Value = -198
$ws.Range("H102").Value = 1551
$ws.Range("I102").Value = 1551
$ws.Range("K102").Value = 1551
$ws.Range("M102").Value = 71
$ws.Range("H110").Value = 2977.182
$ws.Range("J110").Value = 4997.6
$ws.Range("L110").Value = 4997.6
$ws.Range("N110").Value = -9087.6
$ws.Range("H116").Value = 489.6
$ws.Range("I116").Value = 489.6
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 489.6
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1804.4
$ws.Range("N116").ClearContents()
$ws.Range("H122").Value = 2789.9
$ws.Range("I122").Value = 1130.6666
$ws.Range("J122").Value = 3501
$ws.Range("K122").Value = 3391.9998
$ws.Range("L122").Value = 10503
$ws.Range("M122").Value = -941.9998000000001
$ws.Range("N122").Value = -15403

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 489.6
$ws.Range("I3").Value = 489.6
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 489.6
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -375.6
$ws.Range("N3").ClearContents()
$ws.Range("H4").Value = 310
$ws.Range("I4").Value = 310
$ws.Range("K4").Value = 310
$ws.Range("M4").Value = -195
$ws.Range("H86").Value = 3515.818
$ws.Range("I86").Value = 3334.25
$ws.Range("K86").Value = 3334.25
$ws.Range("M86").Value = -2211.25
$ws.Range("H89").Value = 3515.818
$ws.Range("I89").Value = 3334.25
$ws.Range("K89").Value = 16671.25
$ws.Range("M89").Value = -11055.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 765.875
$ws.Range("I22").Value = 406.75
$ws.Range("J22").Value = 1125
$ws.Range("K22").Value = 406.75
$ws.Range("L22").Value = 1125
$ws.Range("M22").Value = -56.75
$ws.Range("N22").Value = -1825

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 1000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 3000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -3540
$ws.Range("H67").Value = 1000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 1000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 3000
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -4872
$ws.Range("H68").Value = 1396.6666
$ws.Range("I68").Value = 1717
$ws.Range("J68").Value = 1305.1428
$ws.Range("K68").Value = 5151
$ws.Range("L68").Value = 3915.4284
$ws.Range("M68").Value = -4340
$ws.Range("N68").Value = -5537.428400000001
$ws.Range("H71").Value = 1396.6666
$ws.Range("I71").Value = 1717
$ws.Range("J71").Value = 1305.1428
$ws.Range("K71").Value = 15453
$ws.Range("L71").Value = 11746.2852
$ws.Range("M71").Value = -11397
$ws.Range("N71").Value = -19858.2852
$ws.Range("H100").Value = 3028
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H114").Value = 500
$ws.Range("I114").Value = 500
$ws.Range("K114").Value = 1500
$ws.Range("M114").Value = 1754
$ws.Range("H117").Value = 272.4
$ws.Range("I117").Value = 308.25
$ws.Range("K117").Value = 924.75
$ws.Range("M117").Value = 2517.25
$ws.Range("H121").Value = 747.4167
$ws.Range("I121").Value = 280.75
$ws.Range("K121").Value = 842.25
$ws.Range("M121").Value = 467.75
$ws.Range("H129").Value = 599
$ws.Range("J129").Value = 913
$ws.Range("L129").Value = 2739
$ws.Range("N129").Value = -12739
$ws.Range("H131").Value = 2000
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("I137").Value = 962.5
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 2887.5
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = 2212.5
$ws.Range("N137").Value = -22200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1008.75
$ws.Range("I102").Value = 938.5714
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 938.5714
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 683.4286
$ws.Range("N102").Value = -4744
$ws.Range("H122").Value = 1243.75
$ws.Range("I122").Value = 1141.6666
$ws.Range("K122").Value = 3424.9998
$ws.Range("M122").Value = -974.9998000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 16734.154
$ws.Range("I7").Value = 15116.111
$ws.Range("J7").Value = 20374.75
$ws.Range("K7").Value = 15116.111
$ws.Range("L7").Value = 20374.75
$ws.Range("M7").Value = -15004.111
$ws.Range("N7").Value = -20598.75
$ws.Range("H22").Value = 1587.125
$ws.Range("I22").Value = 1339.6
$ws.Range("K22").Value = 1339.6
$ws.Range("M22").Value = -1044.6
$ws.Range("H27").Value = 1587.125
$ws.Range("I27").Value = 1339.6
$ws.Range("K27").Value = 1339.6
$ws.Range("M27").Value = -1232.6
$ws.Range("H93").Value = 1800
$ws.Range("I93").Value = 1666.6666
$ws.Range("J93").Value = 2200
$ws.Range("K93").Value = 1666.6666
$ws.Range("L93").Value = 2200
$ws.Range("M93").Value = -418.6666
$ws.Range("N93").Value = -4696
$ws.Range("H100").Value = 3750
$ws.Range("I100").Value = 3750
$ws.Range("K100").Value = 3750
$ws.Range("M100").Value = -3209
$ws.Range("H126").Value = 16734.154
$ws.Range("I126").Value = 15116.111
$ws.Range("J126").Value = 20374.75
$ws.Range("K126").Value = 45348.333
$ws.Range("L126").Value = 61124.25
$ws.Range("M126").Value = -42878.333
$ws.Range("N126").Value = -66064.25
$ws.Range("H136").Value = 7699.6
$ws.Range("I136").Value = 6833
$ws.Range("J136").Value = 8999.5
$ws.Range("K136").Value = 20499
$ws.Range("L136").Value = 26998.5
$ws.Range("M136").Value = -17949
$ws.Range("N136").Value = -32098.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1397.25
$ws.Range("I96").Value = 1429.6666
$ws.Range("J96").Value = 1300
$ws.Range("K96").Value = 1429.6666
$ws.Range("L96").Value = 1300
$ws.Range("M96").Value = -56.66660000000002
$ws.Range("N96").Value = -4046
$ws.Range("H107").Value = 1868.56
$ws.Range("I107").Value = 2306.2144
$ws.Range("J107").Value = 1311.5454
$ws.Range("K107").Value = 6918.6432
$ws.Range("L107").Value = 3934.6362
$ws.Range("M107").Value = -4998.6432
$ws.Range("N107").Value = -7774.6362
$ws.Range("H126").Value = 44375.75
$ws.Range("J126").Value = 53200.4
$ws.Range("L126").Value = 159601.2
$ws.Range("N126").Value = -164541.2
$ws.Range("H136").Value = 5900.143
$ws.Range("I136").Value = 7300.2
$ws.Range("J136").Value = 2400
$ws.Range("K136").Value = 21900.6
$ws.Range("L136").Value = 7200
$ws.Range("M136").Value = -19350.6
$ws.Range("N136").Value = -12300
